$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1. "Curso (semestre ideal): EB (6)" -> "EB (8)"
[void]$d.Content.Find.Execute("Curso (semestre ideal): EB (6)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Curso (semestre ideal): EB (8)", 2)

# 2. Add a new ListBullet paragraph with the docente's name right after the
#    "Docente(s) Responsável(eis) " heading paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Docente(s) Responsável(eis) ") {
        $xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr><w:r><w:t xml:space=`"preserve`">Docente(s) Responsável(eis) </w:t></w:r></w:p>" +
               "<w:p $wns><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr><w:r><w:t>8711290 - Elisson Antônio da Costa Romanel</w:t></w:r></w:p>"
        [void]$p.Range.InsertXML($xml)
        break
    }
}

# 3. Update the requisitos list: replace the whole paragraph so each line
#    keeps its own separate run (avoids Word's run-coalescing on Find/Replace).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("LOT2008 -  Bioquímica II")) {
        $xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr>" +
               "<w:r><w:t>LOT2008 -  Bioquímica II  (Requisito fraco)</w:t><w:br/></w:r>" +
               "<w:r><w:t>LOT2053 -  Microbiologia: da Teoria à Prática  (Requisito fraco)</w:t><w:br/></w:r>" +
               "<w:r><w:t>LOT2040 -  Engenharia Genética Teórica e Prática  (Requisito fraco)</w:t><w:br/></w:r>" +
               "</w:p>"
        [void]$p.Range.InsertXML($xml)
        break
    }
}
